$wb = $excel.ActiveWorkbook
$after = $wb.Worksheets.Item("Terminal Command prompt command")
$newSheet = $wb.Worksheets.Add($null, $after)
$newSheet.Name = "Git_Github"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
